# remove id from template
# The dues table had 3 columns (ID / Unit / Total Due) with a matching
# placeholder row (${ID} / ${UNIT} / ${AMOUNT}). The ID column is removed
# entirely; Word redistributes the freed space across the two remaining
# columns and nudges the table/row metrics, which we replicate explicitly
# below (values taken from the target OOXML, dxa/20 = points).

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Drop the first column ("ID" header / "${ID}" placeholder cell).
$tbl.Columns.Item(1).Delete()

# Resize the remaining two columns / overall table width to match.
$tbl.Columns.Item(1).Width = [double](7473 / 20)
$tbl.Columns.Item(2).Width = [double](5320 / 20)
$tbl.PreferredWidth = [double](12793 / 20)

# Row heights shift slightly once the column is gone.
$tbl.Rows.Item(1).Height = [double](806 / 20)
$tbl.Rows.Item(2).Height = [double](725 / 20)
